# The "Rustic Pizzas" section header (row 38) needs to move down into the
# outlined/grouped block that follows it (rows 39-50 at outlineLevel 1),
# with a new blank separator row taking its old place at row 38.
#
# This is achieved by inserting a new row above the current row 38, which
# shifts the old row 38 (and everything below it) down by one. The newly
# shifted rows are then put at outline level 1 to match the grouping of the
# rest of the pizza section, and the newly inserted blank row 38 has its
# E-column cell's formatting fixed up to match the rest of that (blank) row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row above row 38; old rows 38-49 become rows 39-50.
$ws.Rows(38).Insert()

# The freshly shifted rows (old "Rustic Pizzas" header through the last
# pizza row) belong to the collapsible pizza-section outline group.
$ws.Rows("39:50").OutlineLevel = 1

# Excel's row insert copies formatting from the row above for the new row,
# so E38 picked up row 37's style. Re-copy the correct (blank) formatting
# for that cell from its row-mate D38 so it matches the rest of row 38.
$ws.Range("D38").Copy()
$ws.Range("E38").PasteSpecial(-4122)
